$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws 'D2' '44.256.53'
Set-TextValue $ws 'E2' '  +1.75%  '
Set-TextValue $ws 'D3' '2.361.51'
Set-TextValue $ws 'E3' '  -0.92%  '
Set-TextValue $ws 'E4' '  +0.01%  '
Set-TextValue $ws 'D5' '0.692'
Set-TextValue $ws 'E5' '  +6.19%  '
Set-TextValue $ws 'D6' '243.63'
Set-TextValue $ws 'E6' '  +3.00%  '
Set-TextValue $ws 'D7' '74.17'
Set-TextValue $ws 'E7' '  +2.60%  '
Set-TextValue $ws 'E8' '  +0.01%  '
Set-TextValue $ws 'D9' '0.598'
Set-TextValue $ws 'E9' '  +28.06%  '
Set-TextValue $ws 'E10' '  +5.94%  '
Set-TextValue $ws 'D11' '31.67'
Set-TextValue $ws 'E11' '  +17.41%  '
Set-TextValue $ws 'D12' '7.50'
Set-TextValue $ws 'E12' '  +19.24%  '
Set-TextValue $ws 'E13' '  +2.19%  '
Set-TextValue $ws 'D14' '2.715.87'
Set-TextValue $ws 'E14' '  -0.93%  '
Set-TextValue $ws 'E15' '  +6.64%  '
Set-TextValue $ws 'D16' '0.915'
Set-TextValue $ws 'E16' '  +6.60%  '
Set-TextValue $ws 'D17' '2.356.65'
Set-TextValue $ws 'E17' '  -1.17%  '
Set-TextValue $ws 'D18' '44.291.38'
Set-TextValue $ws 'E18' '  +1.80%  '
Set-TextValue $ws 'E19' '  +4.90%  '
Set-TextValue $ws 'D20' '6.75'
Set-TextValue $ws 'E20' '  +5.44%  '
Set-TextValue $ws 'D21' '78.54'
Set-TextValue $ws 'E21' '  +5.08%  '
Set-TextValue $ws 'D22' '256.87'
Set-TextValue $ws 'E22' '  +1.98%  '
Set-TextValue $ws 'B23' 'Dai'
Set-TextValue $ws 'C23' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws 'D23' '1.00'
Set-TextValue $ws 'E23' '  +0.05%  '
Set-TextValue $ws 'B24' 'PancakeSwap'
Set-TextValue $ws 'C24' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws 'D24' '2.58'
Set-TextValue $ws 'E24' '  +3.43%  '
Set-TextValue $ws 'D25' '3.76'
Set-TextValue $ws 'E25' '  -2.86%  '
Set-TextValue $ws 'D26' '10.83'
Set-TextValue $ws 'E26' '  +7.80%  '
Set-TextValue $ws 'E27' '  +2.55%  '
Set-TextValue $ws 'D28' '22.68'
Set-TextValue $ws 'E28' '  -1.35%  '
Set-TextValue $ws 'E29' '  +5.78%  '
Set-TextValue $ws 'D30' '175.25'
Set-TextValue $ws 'E30' '  +0.36%  '
Set-TextValue $ws 'D31' '0.130'
Set-TextValue $ws 'E31' '  +1.72%  '
Set-TextValue $ws 'D32' '0.135'
Set-TextValue $ws 'E32' '  +5.87%  '
Set-TextValue $ws 'D33' '5.40'
Set-TextValue $ws 'E33' '  +7.61%  '
Set-TextValue $ws 'D34' '0.0758'
Set-TextValue $ws 'E34' '  +9.07%  '
Set-TextValue $ws 'D35' '5.40'
Set-TextValue $ws 'E35' '  +6.25%  '
Set-TextValue $ws 'D36' '3.89'
Set-TextValue $ws 'E36' '  +4.85%  '
Set-TextValue $ws 'E37' '  -0.01%  '
Set-TextValue $ws 'D38' '6.54'
Set-TextValue $ws 'E38' '  -1.24%  '
Set-TextValue $ws 'E39' '  +6.80%  '
Set-TextValue $ws 'D40' '19.16'
Set-TextValue $ws 'E40' '  +0.44%  '
Set-TextValue $ws 'D41' '9.06'
Set-TextValue $ws 'E41' '  +1.59%  '
Set-TextValue $ws 'E42' '  -0.08%  '
Set-TextValue $ws 'E43' '  +16.01%  '
Set-TextValue $ws 'D44' '2.51'
Set-TextValue $ws 'E44' '  +11.37%  '
Set-TextValue $ws 'E45' '  +2.79%  '
Set-TextValue $ws 'D46' '0.1000'
Set-TextValue $ws 'E46' '  +5.03%  '
Set-TextValue $ws 'D47' '101.63'
Set-TextValue $ws 'E47' '  +1.88%  '
Set-TextValue $ws 'E48' '  -0.65%  '
Set-TextValue $ws 'E49' '  -1.75%  '
Set-TextValue $ws 'D50' '1.460.41'
Set-TextValue $ws 'E50' '  +0.59%  '
Set-TextValue $ws 'B51' 'MultiversX'
Set-TextValue $ws 'C51' 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue $ws 'D51' '53.28'
Set-TextValue $ws 'E51' '  +4.55%  '
